$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.340.11'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '3.182.41'
$ws.Range('E3').Value = '  +4.24%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.74'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.22'
$ws.Range('E6').Value = '  +6.90%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.177.68'
$ws.Range('E8').Value = '  +4.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  +4.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +5.58%  '
$ws.Range('E11').Value = '  +3.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.508'
$ws.Range('E12').Value = '  +6.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000275'
$ws.Range('E13').Value = '  +19.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.48'
$ws.Range('E14').Value = '  +10.13%  '
$ws.Range('D15').Value = '3.700.33'
$ws.Range('E15').Value = '  +4.12%  '
$ws.Range('D16').Value = '65.347.23'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.25'
$ws.Range('E17').Value = '  +8.11%  '
$ws.Range('D18').Value = '3.178.97'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '515.68'
$ws.Range('E20').Value = '  +8.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.00'
$ws.Range('E21').Value = '  +7.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.07'
$ws.Range('E22').Value = '  +12.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.742'
$ws.Range('E23').Value = '  +9.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.92'
$ws.Range('E24').Value = '  +4.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.08'
$ws.Range('E25').Value = '  +4.49%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.20'
$ws.Range('E27').Value = '  +16.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.92'
$ws.Range('E28').Value = '  +4.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  +9.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.15'
$ws.Range('E30').Value = '  +7.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.83'
$ws.Range('E31').Value = '  +16.22%  '
$ws.Range('E32').Value = '  +8.14%  '
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.33'
$ws.Range('E34').Value = '  +13.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.67'
$ws.Range('E35').Value = '  +7.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.98'
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '480.32'
$ws.Range('E37').Value = '  +8.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0880'
$ws.Range('E38').Value = '  +9.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.14'
$ws.Range('E39').Value = '  +11.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0423'
$ws.Range('E40').Value = '  +4.39%  '
$ws.Range('D41').Value = '3.133.75'
$ws.Range('E41').Value = '  +5.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.68'
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('E43').Value = '  +6.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  +18.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.292'
$ws.Range('E45').Value = '  +12.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.44'
$ws.Range('E46').Value = '  +6.99%  '
$ws.Range('E47').Value = '  +14.49%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.34'
$ws.Range('E50').Value = '  +13.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.65'
$ws.Range('E51').Value = '  +6.43%  '
